# Working version of PL model
# Adds a "Poland" row to the system_bu_names lookup sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new country/variable-name pair as row 5 (table currently
# runs from row 1 header through row 4 = Hungary).
$ws.Range("A5").Value = "Poland"
$ws.Range("B5").Value = "tu_fa_itjt_pl_HeadID"

Write-Output "Added Poland row to $($ws.Name)"
